$d = $word.ActiveDocument

# 1. Drop the time portion from the test date/time cell.
$d.Content.Find.Execute(
    "10. 12. 2024, 9.30", $true, $false, $false, $false, $false,
    $true, 1, $false, "10. 12. 2024", 2
) | Out-Null

# 2. Fix typo spacing: merge the split "Uporabnik si ogleda rangiranje
#    dobaviteljev." sentence (previously split across three runs around a
#    spell-check proofing mark) back into one continuous run of text.
$d.Content.Find.Execute(
    "Uporabnik si ogleda rangiranje dobaviteljev.", $true, $false, $false, $false, $false,
    $true, 1, $false, "Uporabnik si ogleda rangiranje dobaviteljev.", 2
) | Out-Null
